$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '328.92'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.72%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.08'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '3.13%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.559'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-5.51%'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.69%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.039'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '6.06%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.744'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.91%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4.541'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.00%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9188'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1258'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.23%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1948'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.03%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09318'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2.23%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03734'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '5.08%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.1056'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.94%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001304'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.56%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006257'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '2.24%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.439'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.72%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-2.24%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.314'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-4.74%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1394'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.71%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04418'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.16%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001260'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.11%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004299'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-3.62%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001182'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '3.69%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02753'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '12.86%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05402'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.12%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007668'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.25%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1415'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.57%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008996'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-5.39%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002134'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.65%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01134'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '13.82%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006904'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '2.63%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.18%'
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003555'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '18.46%'
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002284'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '60.53%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002104'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.18%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002004'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.18%'
